$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2083.5557
$ws.Range("I40").Value = 1964.7142
$ws.Range("J40").Value = 2499.5
$ws.Range("K40").Value = 1964.7142
$ws.Range("L40").Value = 2499.5
$ws.Range("M40").Value = -1789.7142
$ws.Range("N40").Value = -2849.5
$ws.Range("H42").Value = 350.5
$ws.Range("J42").Value = 554.5
$ws.Range("L42").Value = 1663.5
$ws.Range("N42").Value = -2123.5
$ws.Range("H88").Value = 2951.5
$ws.Range("J88").Value = 2951.5
$ws.Range("L88").Value = 2951.5
$ws.Range("N88").Value = -3763.5
$ws.Range("H91").Value = 2951.5
$ws.Range("J91").Value = 2951.5
$ws.Range("L91").Value = 2951.5
$ws.Range("N91").Value = -5759.5
$ws.Range("H92").Value = 19231046
$ws.Range("I92").Value = 22727486
$ws.Range("J92").Value = 624
$ws.Range("K92").Value = 22727486
$ws.Range("L92").Value = 624
$ws.Range("M92").Value = -22726238
$ws.Range("N92").Value = -3120
$ws.Range("H98").Value = 2460.3333
$ws.Range("I98").Value = 2624.8333
$ws.Range("K98").Value = 2624.8333
$ws.Range("M98").Value = -1126.8333
$ws.Range("H116").Value = 13351.8
$ws.Range("J116").Value = 3789.7144
$ws.Range("L116").Value = 3789.7144
$ws.Range("N116").Value = -10673.7144
$ws.Range("H122").Value = 2460.3333
$ws.Range("I122").Value = 2624.8333
$ws.Range("K122").Value = 7874.499899999999
$ws.Range("M122").Value = -5424.499899999999
$ws.Range("H132").Value = 1196.28
$ws.Range("J132").Value = 1302.4
$ws.Range("L132").Value = 3907.2
$ws.Range("N132").Value = -8967.200000000001
$ws.Range("H138").Value = 2572.5442
$ws.Range("J138").Value = 2084.7292
$ws.Range("L138").Value = 6254.187600000001
$ws.Range("N138").Value = -16534.1876
$ws.Range("H141").Value = 3441.6
$ws.Range("I141").Value = 1354.8889
$ws.Range("K141").Value = 4064.6667
$ws.Range("M141").Value = 1115.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4349.6895
$ws.Range("I32").Value = 3468.547
$ws.Range("K32").Value = 3468.547
$ws.Range("M32").Value = -3181.547
$ws.Range("H61").Value = 7618.8823
$ws.Range("I61").Value = 8676.916999999999
$ws.Range("J61").Value = 5079.6
$ws.Range("K61").Value = 8676.916999999999
$ws.Range("L61").Value = 5079.6
$ws.Range("M61").Value = -8464.916999999999
$ws.Range("N61").Value = -5503.6
$ws.Range("H63").Value = 1999.6666
$ws.Range("I63").Value = 1999.6666
$ws.Range("K63").Value = 1999.6666
$ws.Range("M63").Value = -1313.6666
$ws.Range("H66").Value = 1999.6666
$ws.Range("I66").Value = 1999.6666
$ws.Range("K66").Value = 9998.333000000001
$ws.Range("M66").Value = -6566.333000000001
$ws.Range("H74").Value = 1453.1428
$ws.Range("J74").Value = 4513.4
$ws.Range("L74").Value = 4513.4
$ws.Range("N74").Value = -6261.4
$ws.Range("H77").Value = 1453.1428
$ws.Range("J77").Value = 4513.4
$ws.Range("L77").Value = 22567
$ws.Range("N77").Value = -31303
$ws.Range("H88").Value = 3972.6
$ws.Range("J88").Value = 3972.6
$ws.Range("L88").Value = 3972.6
$ws.Range("N88").Value = -4784.6
$ws.Range("H91").Value = 3972.6
$ws.Range("J91").Value = 3972.6
$ws.Range("L91").Value = 3972.6
$ws.Range("N91").Value = -6780.6
$ws.Range("H136").Value = 7618.8823
$ws.Range("I136").Value = 8676.916999999999
$ws.Range("J136").Value = 5079.6
$ws.Range("K136").Value = 26030.751
$ws.Range("L136").Value = 15238.8
$ws.Range("M136").Value = -23480.751
$ws.Range("N136").Value = -20338.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31695.223
$ws.Range("J82").Value = 34285.715
$ws.Range("L82").Value = 34285.715
$ws.Range("N82").Value = -35051.715
$ws.Range("H85").Value = 31695.223
$ws.Range("J85").Value = 34285.715
$ws.Range("L85").Value = 34285.715
$ws.Range("N85").Value = -36937.715
$ws.Range("H86").Value = 79611.38
$ws.Range("I86").Value = 3154.8333
$ws.Range("K86").Value = 3154.8333
$ws.Range("M86").Value = -2031.8333
$ws.Range("H89").Value = 79611.38
$ws.Range("I89").Value = 3154.8333
$ws.Range("K89").Value = 15774.1665
$ws.Range("M89").Value = -10158.1665
$ws.Range("H107").Value = 1161.2727
$ws.Range("I107").Value = 996.25
$ws.Range("K107").Value = 996.25
$ws.Range("M107").Value = 923.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -326
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 2302.3333
$ws.Range("I31").Value = 1722.8
$ws.Range("J31").Value = 2716.2856
$ws.Range("K31").Value = 1722.8
$ws.Range("L31").Value = 2716.2856
$ws.Range("M31").Value = -1427.8
$ws.Range("N31").Value = -3306.2856
$ws.Range("H34").Value = 2302.3333
$ws.Range("I34").Value = 1722.8
$ws.Range("J34").Value = 2716.2856
$ws.Range("K34").Value = 1722.8
$ws.Range("L34").Value = 2716.2856
$ws.Range("M34").Value = -1520.8
$ws.Range("N34").Value = -3120.2856
$ws.Range("H132").Value = 2770.9285
$ws.Range("I132").Value = 1307.4445
$ws.Range("J132").Value = 5405.2
$ws.Range("K132").Value = 3922.3335
$ws.Range("L132").Value = 16215.6
$ws.Range("M132").Value = -1392.3335
$ws.Range("N132").Value = -21275.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 663.3333
$ws.Range("J5").Value = 892
$ws.Range("L5").Value = 2676
$ws.Range("N5").Value = -2900
$ws.Range("H135").Value = 663.3333
$ws.Range("J135").Value = 892
$ws.Range("L135").Value = 8028
$ws.Range("N135").Value = -13098

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 808.4286
$ws.Range("I113").Value = 361.75
$ws.Range("J113").Value = 1404
$ws.Range("K113").Value = 361.75
$ws.Range("L113").Value = 1404
$ws.Range("M113").Value = 1808.25
$ws.Range("N113").Value = -5744
$ws.Range("H122").Value = 1032.1666
$ws.Range("I122").Value = 864
$ws.Range("J122").Value = 1200.3334
$ws.Range("K122").Value = 2592
$ws.Range("L122").Value = 3601.0002
$ws.Range("M122").Value = -142
$ws.Range("N122").Value = -8501.0002
$ws.Range("H132").Value = 1604798.4
$ws.Range("I132").Value = 1750007.5
$ws.Range("J132").Value = 7499
$ws.Range("K132").Value = 5250022.5
$ws.Range("L132").Value = 22497
$ws.Range("M132").Value = -5247492.5
$ws.Range("N132").Value = -27557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3207.111
$ws.Range("I7").Value = 3524.1667
$ws.Range("K7").Value = 3524.1667
$ws.Range("M7").Value = -3412.1667
$ws.Range("H82").Value = 1640.6
$ws.Range("I82").Value = 1473.7778
$ws.Range("J82").Value = 1777.091
$ws.Range("K82").Value = 1473.7778
$ws.Range("L82").Value = 1777.091
$ws.Range("M82").Value = -1112.7778
$ws.Range("N82").Value = -2499.091
$ws.Range("H85").Value = 1640.6
$ws.Range("I85").Value = 1473.7778
$ws.Range("J85").Value = 1777.091
$ws.Range("K85").Value = 1473.7778
$ws.Range("L85").Value = 1777.091
$ws.Range("M85").Value = -225.7778000000001
$ws.Range("N85").Value = -4273.091
$ws.Range("H126").Value = 3207.111
$ws.Range("I126").Value = 3524.1667
$ws.Range("K126").Value = 10572.5001
$ws.Range("M126").Value = -8102.500100000001
$ws.Range("H136").Value = 2367.5557
$ws.Range("J136").Value = 2621
$ws.Range("L136").Value = 7863
$ws.Range("N136").Value = -12963

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 79999
$ws.Range("J108").Value = 79999
$ws.Range("L108").Value = 79999
$ws.Range("N108").Value = -87679
$ws.Range("H122").Value = 19990
$ws.Range("I122").Value = 31246.691
$ws.Range("J122").Value = 1697.875
$ws.Range("K122").Value = 93740.073
$ws.Range("L122").Value = 5093.625
$ws.Range("M122").Value = -91290.073
$ws.Range("N122").Value = -9993.625
$ws.Range("H132").Value = 2477.0435
$ws.Range("I132").Value = 1594.1177
$ws.Range("J132").Value = 4978.6665
$ws.Range("K132").Value = 4782.3531
$ws.Range("L132").Value = 14935.9995
$ws.Range("M132").Value = -2252.3531
$ws.Range("N132").Value = -19995.9995
